# This script reproduces the authored edit: a new "Add files via upload" commit
# that appended a new FAQ entry (row 32) about audio-connection sound quality,
# and normalized the "Category" value of two nearby rows (E30/E31) to "All".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize Category ("All") on existing rows 30 and 31 ---
$ws.Range("E30").Value = "All"
$ws.Range("E31").Value = "All"

# --- Append the new FAQ row 32: "차량의 오디오 연결 방식에 따른 음질의 차이" ---
$ws.Range("A32").Value = "차량의 오디오 연결 방식에 따른 음질의 차이"
$ws.Range("B32").Value = 45453
$ws.Range("C32").Value = "HT203025"

$dataLines = @(
    '이 문서는 차량을 유선 CarPlay로 연결하거나, Bluetooth로 연결하였을 때 음질 차이에 대해 설명합니다.',
    '<br>',
    '<br>일반적인 사람의 귀로는 두 연결 방식에 따른 음질의 차이를 느끼기 어렵습니다.',
    '<br>하지만 극히 일부 고객이 유선 CarPlay를 이용해 오디오를 듣는 경우 더 좋다고 느낄 수 있습니다.',
    '<br>',
    '<br>유선 CarPlay는 Bluetooth보다 더 넓은 대역폭을 지원하기 때문에 실제로 음질이 더 좋습니다.',
    '<br>',
    '<br><strong>유선 CarPlay에서 지원하는 코덱은 다음과 같습니다.</strong>',
    '<br> - AAC (Advanced Audio Codec), ALAC (Apple Lossless Audio Codec), Linear PCM (Pulse Code Modulation)',
    '<br>',
    '<br><strong>Bluetooth에서 지원하는 코덱은 다음과 같습니다.</strong>',
    '<br>SBC (Subband Coding), AAC (Advanced Audio Codec), aptX(iPhone 미지원), LDAC',
    '<br>',
    '<br><strong>Bluetooth 4.0/4.1과의 비교</strong>',
    '<br>대역폭: Bluetooth 4.0/4.1의 최대 데이터 전송 속도는 약 25 Mbps로, 유선 또는 무선 CarPlay의 Wi-Fi Direct에 비해 낮습니다. 따라서 고해상도 오디오 스트리밍에는 적합하지 않을 수 있습니다.',
    '<br>손실률: Bluetooth는 데이터 압축을 통해 오디오를 전송하므로, 유선 CarPlay와 비교할 때 오디오 품질이 낮아질 수 있습니다.',
    '<br>',
    '<br><strong>대역폭이란 어떻게 이해하면 됩니까?</strong>',
    '<br>대역폭을 이해할 때에는 차가 다니는 도로를 생각하면 좋습니다. 차량이 소통할 때에는 차로가 많고 폭이 넓은 도로가 소통하기에 적합합니다. 음질에서 말하는 대역폭 또한 넓은 것이 좋은 것입니다.',
    '<br>',
    '<br><strong>고객에게 설명할 때에는 어떻게 하면 좋습니까?</strong>',
    '<br><em>"실제 유선 CarPlay와 Bluetooth 4.0/4.1을 사용할 때 유선 CarPlay가 더 음질이 좋다고 생각하실 수 있습니다. 이는 기술의 차이가 있기 때문입니다. 하지만 일반적인 사람의 귀로는 그 차이를 느끼지 못할 수 있습니다."</em>라고 설명하는 것을 권장합니다.'
)
$ws.Range("D32").Value = [string]::Join("`n", $dataLines)
$ws.Range("D32").WrapText = $true

$ws.Range("E32").Value = "All"

# The new row carries a very long wrapped entry, like the other long rows (4, 21, 29-31)
$ws.Rows.Item(32).RowHeight = 409.6

# Update the view selection to match the newly added row, mirroring the authored workbook
$ws.Range("D32").Select() | Out-Null
